# Update the vm_pu.xlsx results sheet for the 380 kV case (Case_1_2)
# Slack bus voltage setpoint (column B) changes from 1.05 pu to 1.02 pu,
# and all downstream bus voltage results (columns C-F, I-N) are updated
# to the newly recomputed power-flow values for rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020165877422784
$ws.Range("D2").Value = 1.034471823813151
$ws.Range("E2").Value = 1.021209936365017
$ws.Range("F2").Value = 1.032719374532076
$ws.Range("I2").Value = 1.033717567685232
$ws.Range("J2").Value = 1.025364526300644
$ws.Range("K2").Value = 1.037271201087062
$ws.Range("L2").Value = 1.024047940473862
$ws.Range("M2").Value = 1.035523794200081
$ws.Range("N2").Value = 1.012530274150747

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02097605718743
$ws.Range("D3").Value = 1.034978831644508
$ws.Range("E3").Value = 1.021893646607948
$ws.Range("F3").Value = 1.033786552614012
$ws.Range("I3").Value = 1.033872777047784
$ws.Range("J3").Value = 1.025812503273786
$ws.Range("K3").Value = 1.037588209032668
$ws.Range("L3").Value = 1.024538391050495
$ws.Range("M3").Value = 1.036399112437202
$ws.Range("N3").Value = 1.012679966279229

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.021500889404716
$ws.Range("D4").Value = 1.035307139807149
$ws.Range("E4").Value = 1.022336957486079
$ws.Range("F4").Value = 1.034477980019807
$ws.Range("I4").Value = 1.033972224898809
$ws.Range("J4").Value = 1.026102318233085
$ws.Range("K4").Value = 1.037792857941158
$ws.Range("L4").Value = 1.0248559728654
$ws.Range("M4").Value = 1.0369657998536
$ws.Range("N4").Value = 1.012776766674737

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021721669041984
$ws.Range("D5").Value = 1.035445216201757
$ws.Range("E5").Value = 1.022523540266834
$ws.Range("F5").Value = 1.034768868606797
$ws.Range("I5").Value = 1.034013796768048
$ws.Range("J5").Value = 1.026224141934905
$ws.Range("K5").Value = 1.037878777026058
$ws.Range("L5").Value = 1.024989537503933
$ws.Range("M5").Value = 1.037204105591703
$ws.Range("N5").Value = 1.012817446690182

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02175874709114
$ws.Range("D6").Value = 1.035468403044729
$ws.Range("E6").Value = 1.022554880888857
$ws.Range("F6").Value = 1.034817722513041
$ws.Range("I6").Value = 1.034020763008523
$ws.Range("J6").Value = 1.026244595788477
$ws.Range("K6").Value = 1.037893196419656
$ws.Range("L6").Value = 1.025011966681212
$ws.Range("M6").Value = 1.03724412227031
$ws.Range("N6").Value = 1.012824276159938

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02150383892411
$ws.Range("D7").Value = 1.035308984574401
$ws.Range("E7").Value = 1.022339449772083
$ws.Range("F7").Value = 1.034481866053582
$ws.Range("I7").Value = 1.033972781312354
$ws.Range("J7").Value = 1.026103946105704
$ws.Range("K7").Value = 1.037794006450872
$ws.Range("L7").Value = 1.024857757354958
$ws.Range("M7").Value = 1.036968983831348
$ws.Range("N7").Value = 1.012777310302035

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02043955804878
$ws.Range("D8").Value = 1.034643118633423
$ws.Range("E8").Value = 1.021440811121135
$ws.Range("F8").Value = 1.033079847099957
$ws.Range("I8").Value = 1.033770224541198
$ws.Range("J8").Value = 1.025515933183862
$ws.Range("K8").Value = 1.037378433298817
$ws.Range("L8").Value = 1.024213642704889
$ws.Range("M8").Value = 1.035819549799632
$ws.Range("N8").Value = 1.012580875573124

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018568757047853
$ws.Range("D9").Value = 1.033471699431563
$ws.Range("E9").Value = 1.019864297772714
$ws.Range("F9").Value = 1.0306161851209
$ws.Range("I9").Value = 1.033405795169621
$ws.Range("J9").Value = 1.024479397252429
$ws.Range("K9").Value = 1.036642546969375
$ws.Range("L9").Value = 1.023080426629315
$ws.Range("M9").Value = 1.033796426200675
$ws.Range("N9").Value = 1.012234287652907

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017324736253417
$ws.Range("D10").Value = 1.032692156229339
$ws.Range("E10").Value = 1.018818097445319
$ws.Range("F10").Value = 1.028978424172869
$ws.Range("I10").Value = 1.033157841179647
$ws.Range("J10").Value = 1.023788180472377
$ws.Range("K10").Value = 1.036149615120622
$ws.Range("L10").Value = 1.022326229508256
$ws.Range("M10").Value = 1.032449299034644
$ws.Range("N10").Value = 1.012002954083287

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016786834296415
$ws.Range("D11").Value = 1.032354961665597
$ws.Range("E11").Value = 1.018366242863165
$ws.Range("F11").Value = 1.028270377586822
$ws.Range("I11").Value = 1.033049296687816
$ws.Range("J11").Value = 1.023488843692351
$ws.Range("K11").Value = 1.035935630792385
$ws.Range("L11").Value = 1.021999972775257
$ws.Range("M11").Value = 1.031866374673862
$ws.Range("N11").Value = 1.011902723950424

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016587150312988
$ws.Range("D12").Value = 1.032229767503323
$ws.Range("E12").Value = 1.018198579506582
$ws.Range("F12").Value = 1.028007545935096
$ws.Range("I12").Value = 1.033008802002832
$ws.Range("J12").Value = 1.02337765244758
$ws.Range("K12").Value = 1.035856067394311
$ws.Range("L12").Value = 1.021878835142881
$ws.Range("M12").Value = 1.031649910203106
$ws.Range("N12").Value = 1.01186548523199

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016629977898115
$ws.Range("D13").Value = 1.032256619560616
$ws.Range("E13").Value = 1.018234535893965
$ws.Range("F13").Value = 1.028063916582648
$ws.Range("I13").Value = 1.033017496223182
$ws.Range("J13").Value = 1.023401503521496
$ws.Range("K13").Value = 1.035873137617106
$ws.Range("L13").Value = 1.02190481735522
$ws.Range("M13").Value = 1.031696339855583
$ws.Range("N13").Value = 1.011873473454503

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016770325965737
$ws.Range("D14").Value = 1.032344611945873
$ws.Range("E14").Value = 1.01835238016663
$ws.Range("F14").Value = 1.028248648400429
$ws.Range("I14").Value = 1.033045952979716
$ws.Range("J14").Value = 1.023479652669584
$ws.Range("K14").Value = 1.035929055683541
$ws.Range("L14").Value = 1.021989958501326
$ws.Range("M14").Value = 1.031848480432252
$ws.Range("N14").Value = 1.011899645963193

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016856814557617
$ws.Range("D15").Value = 1.032398834291972
$ws.Range("E15").Value = 1.018425011226337
$ws.Range("F15").Value = 1.028362490108893
$ws.Range("I15").Value = 1.033063462768282
$ws.Range("J15").Value = 1.023527802409711
$ws.Range("K15").Value = 1.035963498075097
$ws.Range("L15").Value = 1.022042423254139
$ws.Range("M15").Value = 1.031942227186579
$ws.Range("N15").Value = 1.011915770557722

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017360451326403
$ws.Range("D16").Value = 1.032714542318014
$ws.Range("E16").Value = 1.018848110070527
$ws.Range("F16").Value = 1.029025438455529
$ws.Range("I16").Value = 1.033165020141784
$ws.Range("J16").Value = 1.023808045815513
$ws.Range("K16").Value = 1.036163805254166
$ws.Range("L16").Value = 1.022347888847407
$ws.Range("M16").Value = 1.032487994111477
$ws.Range("N16").Value = 1.012009604772963

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017676575827102
$ws.Range("D17").Value = 1.032912673423332
$ws.Range("E17").Value = 1.01911381979001
$ws.Range("F17").Value = 1.029441587514098
$ws.Range("I17").Value = 1.03322840927963
$ws.Range("J17").Value = 1.023983826393708
$ws.Range("K17").Value = 1.036289308559606
$ws.Range("L17").Value = 1.022539584706715
$ws.Range("M17").Value = 1.032830444497338
$ws.Range("N17").Value = 1.01206844840664

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.017861039829536
$ws.Range("D18").Value = 1.033028273856039
$ws.Range("E18").Value = 1.019268915298693
$ws.Range("F18").Value = 1.029684427711864
$ws.Range("I18").Value = 1.033265269311298
$ws.Range("J18").Value = 1.024086352747881
$ws.Range("K18").Value = 1.036362460173805
$ws.Range("L18").Value = 1.022651428026359
$ws.Range("M18").Value = 1.03303022763873
$ws.Range("N18").Value = 1.012102764973333

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017923949782083
$ws.Range("D19").Value = 1.033067696309915
$ws.Range("E19").Value = 1.019321817723343
$ws.Range("F19").Value = 1.029767248148984
$ws.Range("I19").Value = 1.033277818310112
$ws.Range("J19").Value = 1.024121310961645
$ws.Range("K19").Value = 1.036387394067118
$ws.Range("L19").Value = 1.022689568821598
$ws.Range("M19").Value = 1.033098354881377
$ws.Range("N19").Value = 1.012114465017694

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017642651002366
$ws.Range("D20").Value = 1.032891412309459
$ws.Range("E20").Value = 1.01908530010383
$ws.Range("F20").Value = 1.02939692753348
$ws.Range("I20").Value = 1.033221619989009
$ws.Range("J20").Value = 1.023964967153753
$ws.Range("K20").Value = 1.0362758486473
$ws.Range("L20").Value = 1.022519014405203
$ws.Range("M20").Value = 1.032793698913508
$ws.Range("N20").Value = 1.012062135651998

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016728993723023
$ws.Range("D21").Value = 1.032318698850412
$ws.Range("E21").Value = 1.018317673105584
$ws.Range("F21").Value = 1.028194244829118
$ws.Range("I21").Value = 1.033037578039066
$ws.Range("J21").Value = 1.023456639801611
$ws.Range("K21").Value = 1.035912591407364
$ws.Range("L21").Value = 1.021964885203178
$ws.Range("M21").Value = 1.031803677172075
$ws.Range("N21").Value = 1.01189193905095

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016155217363884
$ws.Range("D22").Value = 1.031958929814074
$ws.Range("E22").Value = 1.017836052094975
$ws.Range("F22").Value = 1.027439045103034
$ws.Range("I22").Value = 1.032920843136309
$ws.Range("J22").Value = 1.023137010239642
$ws.Range("K22").Value = 1.035683734411277
$ws.Range("L22").Value = 1.021616764232641
$ws.Range("M22").Value = 1.031181556442233
$ws.Range("N22").Value = 1.011784879069901

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016459322473189
$ws.Range("D23").Value = 1.03214961934723
$ws.Range("E23").Value = 1.018091271574994
$ws.Range("F23").Value = 1.02783929808012
$ws.Range("I23").Value = 1.032982823017321
$ws.Range("J23").Value = 1.02330645383222
$ws.Range("K23").Value = 1.035805099273166
$ws.Range("L23").Value = 1.021801282600103
$ws.Range("M23").Value = 1.03151132138272
$ws.Range("N23").Value = 1.011841638254549

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017657979942128
$ws.Range("D24").Value = 1.032901019187118
$ws.Range("E24").Value = 1.019098186578026
$ws.Range("F24").Value = 1.029417107124257
$ws.Range("I24").Value = 1.03322468812884
$ws.Range("J24").Value = 1.02397348884339
$ws.Range("K24").Value = 1.036281930764183
$ws.Range("L24").Value = 1.022528309145336
$ws.Range("M24").Value = 1.032810302544939
$ws.Range("N24").Value = 1.012064988132461

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019051849966111
$ws.Range("D25").Value = 1.033774300196419
$ws.Range("E25").Value = 1.020271024053384
$ws.Range("F25").Value = 1.031252280142302
$ws.Range("I25").Value = 1.033500893337683
$ws.Range("J25").Value = 1.024747405311682
$ws.Range("K25").Value = 1.036833209316908
$ws.Range("L25").Value = 1.023373169803165
$ws.Range("M25").Value = 1.034319170541572
$ws.Range("N25").Value = 1.012323939165588

